# Atualização de bases das ligas, do dia: 12-06-2024 às 23:38
# Swap the match-record data (columns B..AD) between each pair of rows
# listed below. Column A (row "id") is left untouched for each row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstCol = 2   # column B
$lastCol  = 30  # column AD

$rowPairs = @(
    @(119, 120),
    @(129, 130),
    @(150, 151)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $cell1 = $ws.Cells.Item($r1, $col)
        $cell2 = $ws.Cells.Item($r2, $col)

        $v1 = $cell1.Value2
        $v2 = $cell2.Value2

        $cell1.Value2 = $v2
        $cell2.Value2 = $v1
    }
}
